$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value2  = "0.4.0-snapshot-1"            # Version
$meta.Range("B6").Value2  = "draft"                        # Status
$meta.Range("B8").Value2  = "2024-05-23T12:16:26+00:00"    # Date
$meta.Range("B10").Value2 = "ANS (https://esante.gouv.fr)" # Contact

# --- Elements sheet: swap the two "Mapping" columns (AK <-> AL) ---
$els = $wb.Worksheets.Item("Elements")

for ($r = 1; $r -le 6; $r++) {
    $akCell = $els.Cells.Item($r, 37)
    $alCell = $els.Cells.Item($r, 38)
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    $akCell.Value2 = $alVal
    $alCell.Value2 = $akVal
}

# Swap the column widths to match the swapped content
# (AK/col 37 now holds the long "business mapping" text -> wide column,
#  AL/col 38 now holds the short "RIM Mapping" text -> narrow column)
$els.Columns.Item(37).ColumnWidth = 86.09375
$els.Columns.Item(38).ColumnWidth = 24.98046875
